$d = $word.ActiveDocument

function Split-ParagraphRuns($paraIndex, $oldLen, $newText, $runLens) {
    # Replace the first $oldLen characters of the paragraph with $newText,
    # then split the resulting single run into separate runs (each keeping
    # its own, empty <w:rPr/>) at the boundaries implied by $runLens.
    $p = $d.Paragraphs.Item($paraIndex)
    $start = $p.Range.Start

    $full = $d.Range($start, $start + $oldLen)
    $full.Text = $newText

    $offset = 0
    foreach ($len in $runLens) {
        $offset = $offset + $len
        $b = $d.Range($start, $start + $offset)
        $b.Bold = 1
        $b.Bold = 0
    }
}

# Locate the three target paragraphs by their distinctive (pre-edit) text,
# disambiguating the two identical "UNIQUE Leg_Number ..." lines by which
# CREATE TABLE block precedes them.
$legScheduleIdx = 0
$legInstanceIdx = 0
$flightDaysIdx = 0
$prevText = ""
$idx = 1
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -eq "UNIQUE Leg_Number DECIMAL(12,0) NOT NULL,`r") {
        if ($prevText -eq "CREATE TABLE Leg_Schedule (`r") {
            $legScheduleIdx = $idx
        } elseif ($prevText -eq "CREATE TABLE Leg_Instance (`r") {
            $legInstanceIdx = $idx
        }
    }
    if ($t -eq "Flight_Number DECIMAL(10,0) NOT NULL,`r") {
        if ($prevText -eq "CREATE TABLE Flight_Days (`r") {
            $flightDaysIdx = $idx
        }
    }
    $prevText = $t
    $idx = $idx + 1
}

if ($legScheduleIdx -eq 0 -or $legInstanceIdx -eq 0 -or $flightDaysIdx -eq 0) {
    Write-Output "ERROR: could not locate target paragraphs (legSchedule=$legScheduleIdx legInstance=$legInstanceIdx flightDays=$flightDaysIdx)"
} else {
    # 1) Leg_Schedule: UNIQUE Leg_Number DECIMAL(12,0) NOT NULL,
    #    -> UNIQUE Leg_Number | INTEGER |  NOT NULL  | AUTO_INCREMENT | ,
    Split-ParagraphRuns $legScheduleIdx 41 `
        "UNIQUE Leg_Number INTEGER NOT NULL AUTO_INCREMENT," `
        @(18, 7, 10, 14, 1)

    # 2) Leg_Instance: UNIQUE Leg_Number DECIMAL(12,0) NOT NULL,
    #    -> UNIQUE Leg_Number | INTEGER  | NOT NULL,
    Split-ParagraphRuns $legInstanceIdx 41 `
        "UNIQUE Leg_Number INTEGER NOT NULL," `
        @(18, 8, 9)

    # 3) Flight_Days: Flight_Number DECIMAL(10,0) NOT NULL,
    #    -> Flight_Number | INTEGER |  NOT NULL,
    Split-ParagraphRuns $flightDaysIdx 37 `
        "Flight_Number INTEGER NOT NULL," `
        @(14, 7, 10)

    Write-Output "OK legScheduleIdx=$legScheduleIdx legInstanceIdx=$legInstanceIdx flightDaysIdx=$flightDaysIdx"
}
